$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B4").Value = 233
$ws.Range("C4").Value = 166

$e4Text = '1: Device root not in the network scenario.
2: Device net not in the network scenario.
3: Devices connected to collision domain A [''as1r1'', ''pc''] are different from the one in the template [''root'', ''as1r2''].
4: Devices connected to collision domain E [''as2r2'', ''dnsnet''] are different from the one in the template [''as2r2'', ''net''].
5: Devices connected to collision domain K [''as1r2'', ''local''] are different from the one in the template [''as3r2'', ''local''].
6: Devices connected to collision domain J [''as3r2'', ''dnsroot''] are different from the one in the template [''as3r2'', ''pc''].
7: The interface `eth0` of `as1r1` has the following IP addresses: [''1.0.0.1/24'']`.
8: The interface `eth1` of `as1r1` has the following IP addresses: [''10.20.0.1/30'']`.
9: The interface `eth0` of `as1r2` has the following IP addresses: [''1.1.0.2/24'']`.
10: The interface `eth1` of `as1r2` has the following IP addresses: [''10.20.1.1/30'']`.
11: The interface `eth2` of `as1r2` has the following IP addresses: [''1.2.0.1/24'']`.
12: The interface `eth2` of `as3r2` has the following IP addresses: [''3.1.0.1/24'']`.
13: Interface eth`3` not found on `as3r2`
14: The interface `eth0` of `local` has the following IP addresses: [''1.2.0.2/24'']`.
15: The interface `eth0` of `pc` has the following IP addresses: [''1.0.0.2/24'']`.
16: `2.1.0.2` not reachable from device `as1r1`.
17: `3.1.0.2` not reachable from device `as1r1`.
18: `3.2.0.1` not reachable from device `as1r1`.
19: `3.2.0.2` not reachable from device `as1r1`.
20: `1.0.0.2` not reachable from device `as1r2`.
21: `2.1.0.2` not reachable from device `as1r2`.
22: `3.1.0.1` not reachable from device `as1r2`.
23: `3.1.0.2` not reachable from device `as1r2`.
24: `3.2.0.1` not reachable from device `as1r2`.
25: `3.2.0.2` not reachable from device `as1r2`.
26: `20.30.0.1` not reachable from device `as1r2`.
27: `1.0.0.2` not reachable from device `as2r1`.
28: `2.1.0.2` not reachable from device `as2r1`.
29: `3.1.0.2` not reachable from device `as2r1`.
30: `3.2.0.1` not reachable from device `as2r1`.
31: `3.2.0.2` not reachable from device `as2r1`.
32: `20.30.1.1` not reachable from device `as2r1`.
33: `1.0.0.2` not reachable from device `as2r2`.
34: `2.1.0.2` not reachable from device `as2r2`.
35: `3.1.0.1` not reachable from device `as2r2`.
36: `3.1.0.2` not reachable from device `as2r2`.
37: `3.2.0.1` not reachable from device `as2r2`.
38: `3.2.0.2` not reachable from device `as2r2`.
39: `1.0.0.2` not reachable from device `as3r1`.
40: `2.1.0.2` not reachable from device `as3r1`.
41: `3.1.0.2` not reachable from device `as3r1`.
42: `3.2.0.1` not reachable from device `as3r1`.
43: `3.2.0.2` not reachable from device `as3r1`.
44: `1.0.0.2` not reachable from device `as3r2`.
45: `2.1.0.2` not reachable from device `as3r2`.
46: `3.1.0.2` not reachable from device `as3r2`.
47: `3.2.0.1` not reachable from device `as3r2`.
48: `3.2.0.2` not reachable from device `as3r2`.
49: `1.0.0.1` not reachable from device `local`.
50: `1.0.0.2` not reachable from device `local`.
51: `1.1.0.1` not reachable from device `local`.
52: `1.1.0.2` not reachable from device `local`.
53: `2.0.0.1` not reachable from device `local`.
54: `2.0.0.2` not reachable from device `local`.
55: `2.1.0.1` not reachable from device `local`.
56: `2.1.0.2` not reachable from device `local`.
57: `3.0.0.1` not reachable from device `local`.
58: `3.0.0.2` not reachable from device `local`.
59: `3.1.0.1` not reachable from device `local`.
60: `3.1.0.2` not reachable from device `local`.
61: `3.2.0.1` not reachable from device `local`.
62: `3.2.0.2` not reachable from device `local`.
63: `10.20.0.1` not reachable from device `local`.
64: `10.20.0.2` not reachable from device `local`.
65: `10.20.1.1` not reachable from device `local`.
66: `10.20.1.2` not reachable from device `local`.
67: `20.30.0.1` not reachable from device `local`.
68: `20.30.0.2` not reachable from device `local`.
69: `20.30.1.1` not reachable from device `local`.
70: `20.30.1.2` not reachable from device `local`.
71: Device `root` is not running.
72: Device `root` is not running.
73: Device `root` is not running.
74: Device `root` is not running.
75: Device `root` is not running.
76: Device `root` is not running.
77: Device `root` is not running.
78: Device `root` is not running.
79: Device `root` is not running.
80: Device `root` is not running.
81: Device `root` is not running.
82: Device `root` is not running.
83: Device `root` is not running.
84: Device `root` is not running.
85: Device `root` is not running.
86: Device `root` is not running.
87: Device `root` is not running.
88: Device `root` is not running.
89: Device `root` is not running.
90: Device `root` is not running.
91: Device `root` is not running.
92: Device `root` is not running.
93: Device `net` is not running.
94: Device `net` is not running.
95: Device `net` is not running.
96: Device `net` is not running.
97: Device `net` is not running.
98: Device `net` is not running.
99: Device `net` is not running.
100: Device `net` is not running.
101: Device `net` is not running.
102: Device `net` is not running.
103: Device `net` is not running.
104: Device `net` is not running.
105: Device `net` is not running.
106: Device `net` is not running.
107: Device `net` is not running.
108: Device `net` is not running.
109: Device `net` is not running.
110: Device `net` is not running.
111: Device `net` is not running.
112: Device `net` is not running.
113: Device `net` is not running.
114: Device `net` is not running.
115: `1.1.0.1` not reachable from device `pc`.
116: `1.1.0.2` not reachable from device `pc`.
117: `2.0.0.1` not reachable from device `pc`.
118: `2.0.0.2` not reachable from device `pc`.
119: `2.1.0.1` not reachable from device `pc`.
120: `2.1.0.2` not reachable from device `pc`.
121: `3.0.0.1` not reachable from device `pc`.
122: `3.0.0.2` not reachable from device `pc`.
123: `3.1.0.1` not reachable from device `pc`.
124: `3.1.0.2` not reachable from device `pc`.
125: `3.2.0.1` not reachable from device `pc`.
126: `3.2.0.2` not reachable from device `pc`.
127: `10.20.0.1` not reachable from device `pc`.
128: `10.20.0.2` not reachable from device `pc`.
129: `10.20.1.1` not reachable from device `pc`.
130: `10.20.1.2` not reachable from device `pc`.
131: `20.30.0.1` not reachable from device `pc`.
132: `20.30.0.2` not reachable from device `pc`.
133: `20.30.1.1` not reachable from device `pc`.
134: `20.30.1.2` not reachable from device `pc`.
135: Device net not in the network scenario.
136: Device net not in the network scenario.
137: Device root not in the network scenario.
138: Device root not in the network scenario.
139: The peering between as1r1 and 1.0.0.2 is not up.
140: The peering between as1r2 and 1.0.0.1 is not up.
141: The route 3.2.0.0/24 IS NOT found in the routing table of `as3r1`.
142: The route 0.0.0.0/0 IS NOT found in the routing table of `root`.
143: The route 1.1.0.0/24 IS NOT found in the routing table of `root`.
144: The route 0.0.0.0/0 IS NOT found in the routing table of `net`.
145: The route 2.1.0.0/24 IS NOT found in the routing table of `net`.
146: The route 0.0.0.0/0 IS NOT found in the routing table of `pc`.
147: The route 3.1.0.0/24 IS NOT found in the routing table of `pc`.
148: The route 0.0.0.0/0 IS NOT found in the routing table of `local`.
149: The route 3.2.0.0/24 IS NOT found in the routing table of `local`.
150: Device `root` is not running.
151: Device `root` is not running.
152: named on local is running but answered with REFUSED when quering for .
153: Device `net` is not running.
154: `resolv.conf` file not found for device `as1r1`
155: `resolv.conf` file not found for device `as1r2`
156: `resolv.conf` file not found for device `as2r1`
157: `resolv.conf` file not found for device `as2r2`
158: `resolv.conf` file not found for device `as3r1`
159: `resolv.conf` file not found for device `as3r2`
160: The local name server for device `pc` has ip `3.2.0.2`
161: `pc.net` not reachable from device `as1r1`.
162: `pc.net` not reachable from device `as1r2`.
163: `pc.net` not reachable from device `as2r1`.
164: `pc.net` not reachable from device `as2r2`.
165: `pc.net` not reachable from device `as3r1`.
166: `pc.net` not reachable from device `as3r2`.
'
$ws.Range("E4").Value = $e4Text
